# Commitment update template: add a new "cleared" column (F) after
# "amount" (E), set page orientation, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in column F, row 1 -- becomes a new shared string "cleared"
$ws.Range("F1").Value = "cleared"

# Page orientation -> portrait (adds <pageSetup orientation="portrait"/>)
$ws.PageSetup.Orientation = 1

# Move / update the saved selection to D10
$ws.Range("D10").Select()
